$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- H1: consolidate the three separate "Note" strings into one wrapped, multi-line note. ---
$note = "Note:" + [char]10 + "Add only numbers in UserID field(Example : Add 123456 if your ID is Q123456/U123456)" + [char]10 + "Valid Values for Follow is ""1"" or ""0""(1 Means Yes and 0 Means No)" + [char]10 + "Valid Values for userRole is ""primary"" or ""secondary"""
$ws.Range("H1").Value = $note
$ws.Range("H1").Interior.Color = 65535
$ws.Range("H1").WrapText = $true
$ws.Range("H1").VerticalAlignment = -4160
$ws.Range("A1").EntireRow.RowHeight = 72.5

# --- H2 / H3 used to hold the now-merged note text: clear them out, keep as blank template cells. ---
$ws.Range("H2").ClearFormats()
$ws.Range("H2").ClearContents()
$ws.Range("H2").NumberFormat = "general"

$ws.Range("H3").ClearFormats()
$ws.Range("H3").ClearContents()
$ws.Range("H3").NumberFormat = "general"

# --- H4: a new blank row added below the note block. ---
$ws.Range("H4").NumberFormat = "general"

# --- Row 2: blank "template" cells under the header row, pre-formatted for data entry. ---
$ws.Range("A2:C2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("D2").NumberFormat = "# ?/?"

# --- Column A now gets an explicit width (previously unset; target ~11.54 chars, bestFit). ---
$ws.Columns(1).ColumnWidth = 10.65

# --- Final selection left on B5. ---
$null = $ws.Range("B5").Select()

Write-Host "done"
